# Update NATMI TPM output data (Proc-Tek.xlsx) with newly recomputed values.
# Sending cluster for the whole sheet moves from MuSCs to ECs, the target
# cluster labels/order shift (MuSCs -> Inflammatory-Mac, Neutrophils -> MuSCs,
# Resolving-Mac -> Neutrophils), and every expression-derived numeric column
# is refreshed with the new TPM-based figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (target cluster: ECs)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Proc"
$ws.Range("C2").Value = "Tek"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.241504
$ws.Range("H2").Value = 0.483008
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 66.82269500000001
$ws.Range("N2").Value = 133.64539
$ws.Range("O2").Value = 0.9484795058090217
$ws.Range("P2").Value = 0.9280825194026978
$ws.Range("Q2").Value = 16.13794813328
$ws.Range("R2").Value = 64.55179253312001
$ws.Range("S2").Value = 0.9484795058090217
$ws.Range("T2").Value = 0.9280825194026978

# Row 3 (target cluster: FAPs)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Proc"
$ws.Range("C3").Value = "Tek"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.241504
$ws.Range("H3").Value = 0.483008
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.804800666666666
$ws.Range("N3").Value = 8.414401999999999
$ws.Range("O3").Value = 0.03981126397570188
$ws.Range("P3").Value = 0.05843268823134937
$ws.Range("Q3").Value = 0.6773705802026666
$ws.Range("R3").Value = 4.064223481216
$ws.Range("S3").Value = 0.03981126397570188
$ws.Range("T3").Value = 0.05843268823134937

# Row 4 (target cluster: Inflammatory-Mac)
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Proc"
$ws.Range("C4").Value = "Tek"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.241504
$ws.Range("H4").Value = 0.483008
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.07387233333333333
$ws.Range("N4").Value = 0.221617
$ws.Range("O4").Value = 0.001048541879565906
$ws.Range("P4").Value = 0.001538989588061868
$ws.Range("Q4").Value = 0.01784046398933333
$ws.Range("R4").Value = 0.107042783936
$ws.Range("S4").Value = 0.001048541879565906
$ws.Range("T4").Value = 0.001538989588061868

# Row 5 (target cluster: MuSCs)
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Proc"
$ws.Range("C5").Value = "Tek"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.241504
$ws.Range("H5").Value = 0.483008
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5329995
$ws.Range("N5").Value = 1.065999
$ws.Range("O5").Value = 0.007565380330087788
$ws.Range("P5").Value = 0.007402687347470469
$ws.Range("Q5").Value = 0.128721511248
$ws.Range("R5").Value = 0.514886044992
$ws.Range("S5").Value = 0.007565380330087788
$ws.Range("T5").Value = 0.007402687347470469

# Row 6 (target cluster: Neutrophils)
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Proc"
$ws.Range("C6").Value = "Tek"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.241504
$ws.Range("H6").Value = 0.483008
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.218072
$ws.Range("N6").Value = 0.6542159999999999
$ws.Range("O6").Value = 0.003095308005622714
$ws.Range("P6").Value = 0.004543115430420422
$ws.Range("Q6").Value = 0.05266526028799999
$ws.Range("R6").Value = 0.315991561728
$ws.Range("S6").Value = 0.003095308005622714
$ws.Range("T6").Value = 0.004543115430420422
